{"js": "const replacements = [\n  [\"Play Moirai Blaze Free Slot | Unique Cluster Pays System\", \"Play Moirai Blaze for Free\"],\n  [\"Cluster Pays system makes for unique gameplay\", \"Impressive graphics and enjoyable soundtrack\"],\n  [\"Special symbols add extra excitement and potential rewards\", \"Unique Cluster Pays system and special symbols\"],\n  [\"Free Spins feature with Vishnu and Ganesha Powerball symbols\", \"Up to 6000 times wager per free spin\"],\n  [\"Visually stunning graphics and enjoyable soundtrack\", \"Medium to high variance and RTP of 96.18%\"],\n  [\"No traditional paylines may be confusing for some players\", \"Lack of traditional paylines\"],\n  [\"May be too volatile for players who prefer lower variance slots\", \"Possible deviation from average RTP\"],\n  [\"Read our review of Moirai Blaze, a visually stunning slot game with a Cluster Pays system and special symbols. Play for free and try your luck today.\", \"Enjoy the visually stunning Moirai Blaze slot game with impressive graphics and rewarding gameplay. Play for free now!\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"Play Moirai Blaze Free Slot | Unique Cluster Pays System\", \"Play Moirai Blaze for Free\"),\n    @(\"Cluster Pays system makes for unique gameplay\", \"Impressive graphics and enjoyable soundtrack\"),\n    @(\"Special symbols add extra excitement and potential rewards\", \"Unique Cluster Pays system and special symbols\"),\n    @(\"Free Spins feature with Vishnu and Ganesha Powerball symbols\", \"Up to 6000 times wager per free spin\"),\n    @(\"Visually stunning graphics and enjoyable soundtrack\", \"Medium to high variance and RTP of 96.18%\"),\n    @(\"No traditional paylines may be confusing for some players\", \"Lack of traditional paylines\"),\n    @(\"May be too volatile for players who prefer lower variance slots\", \"Possible deviation from average RTP\"),\n    @(\"Read our review of Moirai Blaze, a visually stunning slot game with a Cluster Pays system and special symbols. Play for free and try your luck today.\", \"Enjoy the visually stunning Moirai Blaze slot game with impressive graphics and rewarding gameplay. Play for free now!\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
